# Generate Report for Handoff
#
# A new handoff XLIFF was generated for source file
# "0673a674-e86e-4d03-989b-bd9358eb7bf8" (row 5 of both the zh-cn and
# de-de localization-status sheets). Update the "Latest Handoff Datetime"
# (column H) on each locale sheet to the new handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-08-29 06:42:05"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-08-29 06:42:13"
